$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-04-04 13:26:25"

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $newTimestamp
}
